$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

# Remove the three rows describing:
#   405 Atrial Fibrillation or Flutter
#   410 Acute Urinary tract infections UTI
#   411 Sepsis or Septic Shock
# which shifts the remaining rows upward and shrinks the used range.
$ws.Rows.Item(5).Resize(3).EntireRow.Delete()
